$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1
$ws.Range("C2").Value = "String"

$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "date"

$ws.Range("B4").Value = -1
$ws.Range("C4").Value = "date"
$ws.Range("U4").Value = 0

$ws.Range("B5").Value = -1
$ws.Range("C5").Value = "date"
$ws.Range("U5").Value = 0

$ws.Range("B6").Value = -1
$ws.Range("C6").Value = "date"

$ws.Range("B7").Value = -1
$ws.Range("C7").Value = "date"

$ws.Range("B8").Value = -1
$ws.Range("C8").Value = "float"
$ws.Range("E8").Value = 753235
$ws.Range("F8").Value = 20479185.23292093
$ws.Range("G8").Value = 27.1883080750649
$ws.Range("H8").Value = 39.062565
$ws.Range("I8").Value = 35.22311431160598
$ws.Range("J8").Value = 1240.667781808462
$ws.Range("K8").Value = -77.63827000000001
$ws.Range("L8").Value = 77.0429933333333

$ws.Range("B9").Value = -1
$ws.Range("C9").Value = "float"
$ws.Range("E9").Value = 753235
$ws.Range("F9").Value = -49123279.35645875
$ws.Range("G9").Value = -65.21640571197588
$ws.Range("H9").Value = -77.07792000000001
$ws.Range("I9").Value = 35.22439373951801
$ws.Range("J9").Value = 1240.757914316595
$ws.Range("K9").Value = -94.6109883333333
$ws.Range("L9").Value = 77.1871616666667

$ws.Range("B10").Value = -1
$ws.Range("C10").Value = "boolean"

$ws.Range("B11").Value = -1
$ws.Range("C11").Value = "boolean"

$ws.Range("B12").Value = -1
$ws.Range("C12").Value = "boolean"

$ws.Range("B13").Value = -1
$ws.Range("C13").Value = "boolean"

$ws.Range("B14").Value = -1
$ws.Range("C14").Value = "boolean"

$ws.Range("B15").Value = -1
$ws.Range("C15").Value = "boolean"

$ws.Range("B16").Value = -1
$ws.Range("C16").Value = "boolean"

$ws.Range("B17").Value = -1
$ws.Range("C17").Value = "boolean"

$ws.Range("B18").Value = -1
$ws.Range("C18").Value = "boolean"

$ws.Range("B19").Value = -1
$ws.Range("C19").Value = "boolean"

$ws.Range("B20").Value = -1
$ws.Range("C20").Value = "date"
$ws.Range("U20").Value = 0

$ws.Range("B21").Value = -1
$ws.Range("C21").Value = "date"
$ws.Range("U21").Value = 0

$ws.Range("B22").Value = -1
$ws.Range("C22").Value = "integer"
$ws.Range("E22").Value = 821065
$ws.Range("F22").Value = 1645324890
$ws.Range("G22").Value = 2003.891153562751
$ws.Range("H22").Value = 2005
$ws.Range("I22").Value = 87.24449256791843
$ws.Range("J22").Value = 7611.601483433573
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9999
$ws.Range("U22").Value = 0

$ws.Range("B23").Value = -1
$ws.Range("C23").Value = "date"
$ws.Range("U23").Value = 0

$ws.Range("B24").Value = -1
$ws.Range("C24").Value = "date"

$ws.Range("B25").Value = -1
$ws.Range("C25").Value = "date"
$ws.Range("U25").Value = 0

$ws.Range("B26").Value = -1
$ws.Range("C26").Value = "date"
$ws.Range("U26").Value = 0

$ws.Range("B27").Value = -1
$ws.Range("C27").Value = "date"

$ws.Range("B28").Value = -1
$ws.Range("C28").Value = "date"
$ws.Range("U28").Value = 0

$ws.Range("B29").Value = -1
$ws.Range("C29").Value = "boolean"

$ws.Range("B30").Value = -1
$ws.Range("C30").Value = "date"
$ws.Range("U30").Value = 0

$ws.Range("B31").Value = -1
$ws.Range("C31").Value = "date"

$ws.Range("B32").Value = -1
$ws.Range("C32").Value = "date"
$ws.Range("U32").Value = 0

$ws.Range("B33").Value = -1
$ws.Range("C33").Value = "date"
$ws.Range("U33").Value = 0

$ws.Range("B34").Value = -1
$ws.Range("C34").Value = "date"
$ws.Range("U34").Value = 0

$ws.Range("B35").Value = -1
$ws.Range("C35").Value = "date"
$ws.Range("U35").Value = 0

$ws.Range("B36").Value = -1
$ws.Range("C36").Value = "date"
